$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the age-0.5 and age-1.5 sample counts (B1, B2)
$ws.Range("B1").Value = 39
$ws.Range("B2").Value = 62

# Add the two total rows beneath the existing histogram data.
# Write the shared-string labels in this order ("Total" then "Total (5.5+)")
# so the generated shared-strings table matches the expected index order.
$ws.Range("A19").Value = "Total"
$ws.Range("A18").Value = "Total (5.5+)"

$ws.Range("B18").Formula = "=SUM(B6:B17)"
$ws.Range("B19").Formula = "=SUM(B1:B17)"

# Move the active selection to A20, matching the post-edit workbook state.
$ws.Range("A20").Select()
